$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tmp = $ws.Range("Z1")
$tmp.NumberFormat = "@"
$tmp.Value = "1"

$tmp.Copy()
$ws.Range("B11").PasteSpecial(-4163)

$tmp.Clear()
